$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 3 ("Main Menu Screen" -> "Menu Screen")
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$title3 = $s3.Shapes.Item(1)
$title3.TextFrame.TextRange.Text = "Menu Screen"

# ---------------------------------------------------------------------------
# Slide 5 ("Find a Band" title -> "Band Search Screen  Find a Band")
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$title5 = $s5.Shapes.Item(1)
$tr5 = $title5.TextFrame.TextRange

# Prepend "Band Search Screen " (keeps the original run formatting).
[void]$tr5.InsertBefore("Band Search Screen ")

# Locate "Find a Band" inside the updated text and insert a single space
# right before it - this becomes its own run.
$full5 = $tr5.Text
$idx5 = $full5.IndexOf("Find a Band")
$findRange5 = $tr5.Characters($idx5 + 1, 11)
[void]$findRange5.InsertBefore(" ")

# Give that lone space run the Wingdings font, matching the authored markup.
$spaceRange5 = $tr5.Characters($idx5 + 1, 1)
$spaceRange5.Font.Name = "Wingdings"

# Resize/reposition the border textbox (TextBox 7) under the new, taller title.
$box5 = $s5.Shapes.Item(2)
$box5.Left = 360.6990661621094
$box5.Top = 133.12496948242188
$box5.Width = 237.8009490966797
$box5.Height = 391.1249694824219

# ---------------------------------------------------------------------------
# Slide 6 ("Find a Musician" title -> "Musician Search Screen  Find a Musician")
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$title6 = $s6.Shapes.Item(1)
$tr6 = $title6.TextFrame.TextRange

[void]$tr6.InsertBefore("Musician Search Screen ")

$full6 = $tr6.Text
$idx6 = $full6.IndexOf("Find a Musician")
$findRange6 = $tr6.Characters($idx6 + 1, 16)
[void]$findRange6.InsertBefore(" ")

$spaceRange6 = $tr6.Characters($idx6 + 1, 1)
$spaceRange6.Font.Name = "Wingdings"

# Resize/reposition the border textbox (TextBox 4) under the new, taller title.
$box6 = $s6.Shapes.Item(2)
$box6.Left = 362.79608154296875
$box6.Top = 133.12496948242188
$box6.Width = 235.70387268066406
$box6.Height = 391.1249694824219
